$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2025-06-13T15:45:04+00:00"
$meta.Range("B15").Value = "4.0.1"

# --- Elements sheet updates ---
$elements = $wb.Worksheets.Item("Elements")

# AJ2 previously held the "long" ele-1/ext-1 constraint text; it now matches the
# shorter text used elsewhere on the sheet (e.g. AJ4), which causes the shared
# string table to de-duplicate those two entries when the workbook is saved.
$elements.Range("AJ2").Value = "ele-1:All FHIR elements must have a @value or children {hasValue() or (children().count() > id.count())}" + [char]10 + "ext-1:Must have either extensions or value[x], not both {extension.exists() != value.exists()}"

# K3 type changed from "id" to "string"
$elements.Range("K3").Value = "string" + [char]10

# M6 documentation link updated from R4B to R4
$elements.Range("M6").Value = "Value of extension - must be one of a constrained set of the data types (see [Extensibility](http://hl7.org/fhir/R4/extensibility.html) for a list)."
